$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Change the Period value in C10 from 4 (quarter) to 2 (8-week period)
$ws.Range("C10").Value = 2

# Recalculate the workbook so all dependent formulas update
$excel.Calculate()

# Update the active cell selection on the frozen bottom-right pane to B10
$ws.Range("B10").Select()
